$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.533.49'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.662.04'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.80'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.32'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '2.661.44'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.170'
$ws.Range('E11').Value = '  +2.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').Value = '3.149.41'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000186'
$ws.Range('E15').Value = '  -2.07%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '72.444.01'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.26'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '2.660.61'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.34'
$ws.Range('E19').Value = '  +5.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.83'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '371.60'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.20'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.06'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.10'
$ws.Range('E24').Value = '  -2.31%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.32'
$ws.Range('E26').Value = '  -2.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').Value = '2.797.71'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = '0.0₃0974'
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.15'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '496.02'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.30'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.48'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.55'
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.112'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.96'
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').Value = '  -5.06%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.01'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.41'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '155.78'
$ws.Range('E47').Value = '  +3.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.75'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.559'
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.73'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0758'
$ws.Range('E51').Value = '  -1.08%  '
